# Reorders the comma-separated email lists in column G ("Recorded By")
# for the affected rows, per the author's sync of attendance_reports.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg"
    3  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    4  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    5  = "asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    6  = "mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    7  = "menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
    11 = "Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    12 = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
    13 = "amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
    20 = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    24 = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
    27 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    28 = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    30 = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
